$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column A: a running "rank" counter next to the keyword list, plus a
# "Total" header that counts how many keywords are present.
$ws.Range("A1").Formula = '="Total: " & COUNTA(B2:B10000)'

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Formula = '=IF(B{0}<>"", COUNTA($B$2:B{0}), "")' -f $r
}

# Match the author's final view state: zoomed in, selection moved to F14.
$ws.Range("F14").Select()
$excel.ActiveWindow.Zoom = 177
